$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: add new value in A2, keep B2/C2 the same values
$ws.Range("A2").Value = "crossover1"
$ws.Range("B2").Value = "环太平洋2"
$ws.Range("C2").Value = "2023.12.21-2024.01.18"

# Update header E1: "特殊联动物品" -> "特殊联动物品/道具"
$ws.Range("E1").Value = "特殊联动物品/道具"

# Update selection/view state to match the new active cell
$ws.Range("E1").Select()
$excel.ActiveWindow.ScrollColumn = 2
